$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1475
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 80
$ws.Range("G2").Value = 53
$ws.Range("H2").Value = 33
$ws.Range("I2").Value = 33
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 2661
$ws.Range("L2").Value = 1263
$ws.Range("M2").Value = 1398
$ws.Range("N2").Value = 1398
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 152
$ws.Range("Q2").Value = -10
$ws.Range("R2").Value = -67
$ws.Range("S2").Value = 42
$ws.Range("T2").Value = 105
$ws.Range("U2").Value = -115
$ws.Range("V2").Value = 680
$ws.Range("W2").Value = 5.44
$ws.Range("X2").Value = 2.21
$ws.Range("Y2").Value = 2.34
$ws.Range("Z2").Value = 1.25
$ws.Range("AA2").Value = 90.38
$ws.Range("AB2").Value = 821.09
$ws.Range("AC2").Value = 216
$ws.Range("AD2").Value = 16.52
$ws.Range("AE2").Value = 9215
$ws.Range("AF2").Value = 0.39
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 2.8
$ws.Range("AI2").Value = 46.57
$ws.Range("AJ2").Value = 15167224

# Row 3
$ws.Range("D3").Value = 1809
$ws.Range("E3").Value = 121
$ws.Range("F3").Value = 121
$ws.Range("G3").Value = 107
$ws.Range("H3").Value = 79
$ws.Range("I3").Value = 82
$ws.Range("J3").Value = -3
$ws.Range("K3").Value = 2817
$ws.Range("L3").Value = 1359
$ws.Range("M3").Value = 1458
$ws.Range("N3").Value = 1455
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 155
$ws.Range("Q3").Value = -10
$ws.Range("R3").Value = 40
$ws.Range("S3").Value = -25
$ws.Range("T3").Value = 36
$ws.Range("U3").Value = -46
$ws.Range("V3").Value = 661
$ws.Range("W3").Value = 6.7
$ws.Range("X3").Value = 4.36
$ws.Range("Y3").Value = 5.73
$ws.Range("Z3").Value = 2.88
$ws.Range("AA3").Value = 93.2
$ws.Range("AB3").Value = 838.13
$ws.Range("AC3").Value = 530
$ws.Range("AD3").Value = 23.32
$ws.Range("AE3").Value = 9389
$ws.Range("AF3").Value = 1.32
$ws.Range("AG3").Value = 130
$ws.Range("AH3").Value = 1.05
$ws.Range("AI3").Value = 24.64
$ws.Range("AJ3").Value = 15501672

# Row 4
$ws.Range("D4").Value = 2091
$ws.Range("E4").Value = 144
$ws.Range("F4").Value = 144
$ws.Range("G4").Value = 107
$ws.Range("H4").Value = 74
$ws.Range("I4").Value = 77
$ws.Range("J4").Value = -4
$ws.Range("K4").Value = 3005
$ws.Range("L4").Value = 1461
$ws.Range("M4").Value = 1544
$ws.Range("N4").Value = 1545
$ws.Range("O4").Value = -1
$ws.Range("P4").Value = 167
$ws.Range("Q4").Value = 109
$ws.Range("R4").Value = -55
$ws.Range("S4").Value = 55
$ws.Range("T4").Value = 50
$ws.Range("U4").Value = 59
$ws.Range("V4").Value = 701
$ws.Range("W4").Value = 6.86
$ws.Range("X4").Value = 3.52
$ws.Range("Y4").Value = 5.17
$ws.Range("Z4").Value = 2.53
$ws.Range("AA4").Value = 94.67
$ws.Range("AB4").Value = 826.05
$ws.Range("AC4").Value = 499
$ws.Range("AD4").Value = 16.7
$ws.Range("AE4").Value = 9266
$ws.Range("AF4").Value = 0.9
$ws.Range("AG4").Value = 130
$ws.Range("AH4").Value = 1.56
$ws.Range("AI4").Value = 27.98
$ws.Range("AJ4").Value = 16672240

# Row 5
$ws.Range("D5").Value = 2329
$ws.Range("E5").Value = 121
$ws.Range("F5").Value = 121
$ws.Range("G5").Value = 99
$ws.Range("H5").Value = 67
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = -3
$ws.Range("K5").Value = 3047
$ws.Range("L5").Value = 1454
$ws.Range("M5").Value = 1592
$ws.Range("N5").Value = 1596
$ws.Range("O5").Value = -4
$ws.Range("P5").Value = 167
$ws.Range("Q5").Value = 121
$ws.Range("R5").Value = -61
$ws.Range("S5").Value = -45
$ws.Range("T5").Value = 46
$ws.Range("U5").Value = 76
$ws.Range("V5").Value = 678
$ws.Range("W5").Value = 5.2
$ws.Range("X5").Value = 2.89
$ws.Range("Y5").Value = 4.48
$ws.Range("Z5").Value = 2.23
$ws.Range("AA5").Value = 91.34
$ws.Range("AB5").Value = 857.4
$ws.Range("AC5").Value = 422
$ws.Range("AD5").Value = 15.11
$ws.Range("AE5").Value = 9576
$ws.Range("AF5").Value = 0.67
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 2.35
$ws.Range("AI5").Value = 35.52
$ws.Range("AJ5").Value = 16672240

# Row 6
$ws.Range("D6").Value = 2472
$ws.Range("E6").Value = 74
$ws.Range("F6").Value = 74
$ws.Range("G6").Value = 56
$ws.Range("H6").Value = 39
$ws.Range("I6").Value = 42
$ws.Range("K6").Value = 3017
$ws.Range("L6").Value = 1416
$ws.Range("M6").Value = 1601
$ws.Range("N6").Value = 1608
$ws.Range("P6").Value = 167
$ws.Range("Q6").Value = 130
$ws.Range("R6").Value = -79
$ws.Range("S6").Value = -50
$ws.Range("T6").Value = 67
$ws.Range("U6").Value = 63
$ws.Range("V6").Value = 655
$ws.Range("W6").Value = 2.98
$ws.Range("X6").Value = 1.59
$ws.Range("Y6").Value = 2.62
$ws.Range("Z6").Value = 1.29
$ws.Range("AA6").Value = 88.44
$ws.Range("AB6").Value = 865.1
$ws.Range("AC6").Value = 252
$ws.Range("AD6").Value = 18.15
$ws.Range("AE6").Value = 9677
$ws.Range("AF6").Value = 0.47
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 3.29
$ws.Range("AI6").Value = 59.45
$ws.Range("AJ6").Value = 16672240

# Row 7: clear all forecast data, keep only A/B/C
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all forecast data, keep only A/B/C
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all forecast data, keep only A/B/C
$ws.Range("D9:AJ9").ClearContents()
